# Create the new "SystemStructureCreateOrgsApps" worksheet, as the last tab.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "SystemStructureCreateOrgsApps"

# Header row
$ws.Range("A1").Value = "#"
$ws.Range("B1").Value = "Организация"
$ws.Range("C1").Value = "Приложение"
$ws.Range("A1:C1").Font.Bold = $true

# Data rows: Org/App pairs
$data = @(
    @(1, "АдмГор", "Тестовое имя 4UFDZEKJ"),
    @(2, "АдмГор", "Тестовое имя 4VX8YZPI"),
    @(3, "АдмГор", "Тестовое имя 4XBJSER7"),
    @(4, "АдмГор", "Тестовое имя 4XRVE1WJ"),
    @(5, "АдмГор", "Тестовое имя 4ZUS6A3J")
)

$row = 2
foreach ($d in $data) {
    $ws.Cells.Item($row, 1).Value = $d[0]
    $ws.Cells.Item($row, 2).Value = $d[1]
    $ws.Cells.Item($row, 3).Value = $d[2]
    $row++
}

# Column widths to fit the new content
$ws.Columns.Item(2).ColumnWidth = 16.28515625
$ws.Columns.Item(3).ColumnWidth = 27.7109375

# Selection on the newly active sheet
$ws.Range("C2").Select()
